$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Salesforce -> NetSuite
$ws.Range("A2").Value = 'NetSuite'
$ws.Range("B2").Value = 'This system is primarily used for financial management, accounting, and reporting. It supports deal entries, purchase order creation, workflow approvals, payment processing integration, journal entries, financial reporting, and general ledger management.'
$ws.Range("C2").Value = 'The client utilizes this system to manage core financial processes such as accounting entries, purchase order workflows, bill approvals, and financial reporting. It also supports integrations for payment processing and management reporting.'
$ws.Range("D2").Value = 'The NetSuite system administration is managed by Blackforge Consulting, led by Steve Monti and Andrew Baker.'
$ws.Range("E2").Value = 'Access provisioning is initiated when Jade team members send an email to Blackforge Consulting''s support email, prompting Blackforge to create a ticket in their ticketing system. The ticket requires approval from authorized approvers, Megan Hodgson or John, before Blackforge executes the changes and documents them in the system.'
$ws.Range("F2").Value = 'Upon termination, HR initiates the removal process by notifying Blackforge, who manually removes access immediately. The removal is documented in a case and shared with HR for record-keeping.'
$ws.Range("G2").Value = 'Access is configured using a role-based model, where custom roles are utilized and permissions are finalized during the implementation process, with global permissions disabled in the system.'
$ws.Range("H2").Value = 'Yes; System administrators can modify roles in NetSuite, and changes have been made since the initial rollout, though these have been minimal.'
$ws.Range("I2").Value = 'Yes; management performs an annual review of all roles and permissions.'
$ws.Range("J2").Value = 'Yes; Users with privileged access include Steve Monti, Andrew Baker, and Justin, who hold the Administrator role in NetSuite, granting them the ability to manage system configurations and integrations.'
$ws.Range("K2").Value = 'Yes; All user accounts are tied to individual employees; there are no generic or shared credentials that can be used interactively.'
$ws.Range("L2").Value = 'Credentials for shared and generic accounts are currently stored securely by Steve Monti until a formal enterprise credential management solution is implemented.'
$ws.Range("M2").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript'
$ws.Range("N2").Value = 'Yes; Reviews are conducted quarterly. Management performs periodic user access reviews for the NetSuite system, starting from its initial implementation, with plans to continue on an ongoing basis.'
$ws.Range("O2").Value = 'The system maintains logs of administrative activities. Logs are retained indefinitely and include detailed audit trails.'
$ws.Range("P2").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("Q2").Value = 'Users authenticate through Azure SSO for the production environment, while the sandbox environment requires independent login credentials without SSO integration.'
$ws.Range("R2").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("S2").Value = 'Management can perform configuration changes to adjust system settings, update workflows, modify integrations, and implement custom scripts. These changes are managed through a ticketing system, developed and tested in a sandbox environment, and pushed to production following approval processes.'
$ws.Range("T2").Value = 'The following individuals have change capabilities: Steve Monti and Andrew Baker from Blackforge Consulting (application changes) and Justin from Whiplate (integration setup and related changes). These individuals have administrator access to NetSuite and are responsible for their respective areas of modification.'
$ws.Range("U2").Value = 'Management maintains separate Sandbox and Release Preview environments for implementing and testing changes.'
$ws.Range("V2").Value = 'The client''s change management process includes the following steps: Change requests are submitted via NetSuite''s ticketing system or email depending on the change type, development is performed in the sandbox environment or directly in the live environment by the responsible team, testing is conducted by Jade''s team with relevant business users or in the live environment for integration changes, approvals are documented either through sign-off or email, and deployment is handled by the Blackforge team or integration owner following the completion of testing and approval.'
$ws.Range("W2").Value = 'Yes, the vendor provides two major system upgrades per year along with occasional patches and bug fixes. The client receives a release preview environment for testing prior to upgrades, and testing evidence is documented. Updates are deployed to Production by the vendor without requiring explicit approval from the client, but notifications are typically sent.'
$ws.Range("X2").Value = 'No, there is no inherent system functionality preventing users from both developing and deploying changes.'
$ws.Range("Y2").Value = 'Management performs periodic reviews of system changes to ensure they were appropriately requested, tested, and approved prior to deployment. The process involves creating a ticket in the ticketing system, testing changes in a sandbox environment, obtaining sign-off from relevant personnel, and documenting all approvals within the ticketing system.'
$ws.Range("Z2").Value = 'Yes, the system has automated jobs and interfaces with Bill.com, Concur, and Adaptive Planning. These jobs perform data integration and exchange functions.'
$ws.Range("AA2").Value = 'Jobs are managed using NetSuite and Adaptive. These tools provide integration management and monitoring capabilities for scheduling and execution.'
$ws.Range("AB2").Value = 'Job failures are handled through a defined process: detection occurs in Adaptive, notifications are sent for both successful and failed tasks, and resolution involves creating missing accounts in Adaptive, remapping, and rerunning the task.'
$ws.Range("AC2").Value = 'Data is stored in a vendor-managed system managed by NetSuite.'
$ws.Range("AD2").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("AE2").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("AF2").Value = 'N/A - Backup failure resolution is managed by the vendor as part of their SaaS service.'
$ws.Range("AG2").Value = 'No, management does not perform regular SOC report reviews.'

# Row 3: Concur -> Freshworks
$ws.Range("A3").Value = 'Freshworks'
$ws.Range("B3").Value = 'This is an IT Service Management tool used for ticket intake, change control, service requests, incidents, problem management, and procurement activities. It supports intake management for IT, security, and HR teams and is currently in phase one of its rollout.'
$ws.Range("C3").Value = 'The client utilizes this system to manage IT Service Management processes, including ticket intake, change control, service requests, incidents, problem management, procurement activities, and select project tasks, while also supporting intake management for HR and Security departments.'
$ws.Range("D3").Value = 'The IT Administration Team, overseen by Andy Masterton (Head of First-Line Support).'
$ws.Range("E3").Value = 'Access provisioning is managed through Freshservice, where HR triggers an automated process for new hires via a lever email, generating tasks for agents. Additional access requests and role changes are handled manually, with approvals documented in a SharePoint spreadsheet and overseen by designated approvers.'
$ws.Range("F3").Value = 'Access removal for terminations and role changes is triggered by an automated process initiated by a lever email from HR, which generates tasks assigned to an agent for access revocation.'
$ws.Range("G3").Value = 'Access is configured using a combination of group-based and individual assignments, where user portal access is provisioned via Active Directory groups, agent access is granted individually, and administrative access is restricted to specific IT personnel with elevated permissions.'
$ws.Range("H3").Value = 'Yes; System administrators can modify roles within Freshservice, but only designated administrators with elevated access are permitted to make these changes.'
$ws.Range("I3").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("J3").Value = 'Yes; Users with privileged access include admins who can access back-end configurations and agents who can perform activities beyond end-user capabilities, with admin access restricted to specific IT personnel approved by Danielle Corfe and Mark Bretner.'
$ws.Range("K3").Value = 'No; All user accounts are tied to individual employees; there are no generic or shared credentials that can be used interactively.'
$ws.Range("L3").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("M3").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("N3").Value = 'No; there is no periodic review process currently in place. Admin access requires approval, and user access is managed via Active Directory groups, but no reviews are conducted to validate ongoing appropriateness.'
$ws.Range("O3").Value = 'Activity logging captures changes related to change control, approvals, and sandbox testing and includes timestamps, documented actions, and compliance reviews. Logs are stored within the system for traceability purposes.'
$ws.Range("P3").Value = 'No; management does not perform structured periodic reviews of user activity, and monitoring is reactive or event-driven rather than proactive.'
$ws.Range("Q3").Value = 'Authentication is managed through Active Directory groups for requester-level access, while elevated agent access is assigned individually, and administrative access requires approval from designated approvers.'
$ws.Range("R3").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("S3").Value = 'Management can perform configuration changes to adjust system settings, update workflows, and modify code within the system. These changes are managed through Freshservice, with approvals tracked and testing conducted in a sandbox environment prior to deployment.'
$ws.Range("T3").Value = 'Only elevated users within the IT department are designated as administrators, with Danielle Corfe and Mark Bretner identified as approvers for granting admin access to Freshservice. Administrative access is restricted to this subset of the IT team, and not all IT personnel have the ability to make configurations or changes.'
$ws.Range("U3").Value = 'Management maintains a sandbox environment for developing and testing changes before deployment to the production environment.'
$ws.Range("V3").Value = 'The client''s change management process includes the following steps: Change requests for Freshservice configuration changes are initiated via change control tickets containing impacted parties, technical steps, rollback plans, and risk assessments. Changes are developed and tested in a sandbox environment before approval by designated change approvers, with documentation maintained in the ticket. Approved changes are deployed to production following testing. Code changes are managed in GitLab, where internal teams create issues, develop changes, conduct testing, and follow approval workflows involving review and merge requests. Production deployments occur after testing and approval, adhering to the software development lifecycle methodology.'
$ws.Range("W3").Value = 'Yes, the vendor pushes updates, patches, and bug fixes to the SaaS system. Changes are first tested in a sandbox environment before being deployed to the production environment.'
$ws.Range("X3").Value = 'The system enforces segregation of duties through role-based access controls and a sandbox environment. Elevated users with administrative roles can make changes, but deployment to production requires a manual approval process involving designated approvers.'
$ws.Range("Y3").Value = 'No periodic review of changes is performed.'
$ws.Range("Z3").Value = 'Yes, the system has two types of automated jobs/interfaces with other in-scope systems. These jobs perform termination process automation and change testing in a sandbox environment.'
$ws.Range("AA3").Value = 'Jobs are managed using Freshservice, Workday, and GitLab. These tools provide capabilities for scheduling, monitoring, and managing workflows, approvals, and automation.'
$ws.Range("AB3").Value = 'Job failures were not discussed in the walkthrough meeting transcript, and no defined process for detection, notification, or resolution was provided.'
$ws.Range("AC3").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("AD3").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("AE3").Value = 'Not applicable - backup implementation is handled by the SaaS vendor.'
$ws.Range("AF3").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("AG3").Value = 'No, management does not perform regular SOC report reviews.'

# Remove the old rows 4 (Concur) and 5 (Microsoft Azure); dimension becomes A1:AG3
$ws.Rows("4:5").Delete()

